$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Gnas", "Gcgr", "ECs", 3, 1, 175.411433, 526.234299, 0.1535106429347505, 0.1535106429347505, 1, 0.3333333333333333, 0.3424073333333333, 1.027222, 0.1616399685287176, 0.1616399685287176, 60.06216100970867, 540.559449087378, 0.02481345549279627, 0.02481345549279627),
    @("ECs", "Gnas", "Gcgr", "FAPs", 3, 1, 175.411433, 526.234299, 0.1535106429347505, 0.1535106429347505, 3, 1, 1.412430333333333, 4.237291, 0.666764909520063, 0.666764909520063, 247.7564287826676, 2229.807859044009, 0.1023555099467556, 0.1023555099467556),
    @("ECs", "Gnas", "Gcgr", "MuSCs", 3, 1, 175.411433, 526.234299, 0.1535106429347505, 0.1535106429347505, 2, 0.6666666666666666, 0.231231, 0.6936929999999999, 0.1091570416994492, 0.1091570416994492, 40.56056106402299, 365.0450495762069, 0.01675676765213782, 0.01675676765213783),
    @("ECs", "Gnas", "Gcgr", "Resolving-Mac", 3, 1, 175.411433, 526.234299, 0.1535106429347505, 0.1535106429347505, 1, 0.3333333333333333, 0.1322646666666667, 0.396794, 0.06243808025177026, 0.06243808025177026, 23.20073471526733, 208.806612437406, 0.009584909843060801, 0.009584909843060804),
    @("FAPs", "Gnas", "Gcgr", "ECs", 3, 1, 531.1103823333333, 1593.331147, 0.464799214434963, 0.4647992144349631, 1, 0.3333333333333333, 0.3424073333333333, 1.027222, 0.1616399685287176, 0.1616399685287176, 181.8560897204038, 1636.704807483634, 0.07513013039344008, 0.07513013039344008),
    @("FAPs", "Gnas", "Gcgr", "FAPs", 3, 1, 531.1103823333333, 1593.331147, 0.464799214434963, 0.4647992144349631, 3, 1, 1.412430333333333, 4.237291, 0.666764909520063, 0.666764909520063, 750.1564143558639, 6751.407729202776, 0.3099118061577245, 0.3099118061577245),
    @("FAPs", "Gnas", "Gcgr", "MuSCs", 3, 1, 531.1103823333333, 1593.331147, 0.464799214434963, 0.4647992144349631, 2, 0.6666666666666666, 0.231231, 0.6936929999999999, 0.1091570416994492, 0.1091570416994492, 122.809184817319, 1105.282663355871, 0.05073610723194851, 0.05073610723194852),
    @("FAPs", "Gnas", "Gcgr", "Resolving-Mac", 3, 1, 531.1103823333333, 1593.331147, 0.464799214434963, 0.4647992144349631, 1, 0.3333333333333333, 0.1322646666666667, 0.396794, 0.06243808025177026, 0.06243808025177026, 70.24713768252421, 632.2242391427179, 0.02902117065185, 0.02902117065185),
    @("MuSCs", "Gnas", "Gcgr", "ECs", 3, 1, 360.115397, 1080.346191, 0.3151536087398187, 0.3151536087398188, 1, 0.3333333333333333, 0.3424073333333333, 1.027222, 0.1616399685287176, 0.1616399685287176, 123.3061527790447, 1109.755375011402, 0.05094141939841607, 0.05094141939841608),
    @("MuSCs", "Gnas", "Gcgr", "FAPs", 3, 1, 360.115397, 1080.346191, 0.3151536087398187, 0.3151536087398188, 3, 1, 1.412430333333333, 4.237291, 0.666764909520063, 0.666764909520063, 508.6379102231757, 4577.741192008581, 0.2101333674163266, 0.2101333674163266),
    @("MuSCs", "Gnas", "Gcgr", "MuSCs", 3, 1, 360.115397, 1080.346191, 0.3151536087398187, 0.3151536087398188, 2, 0.6666666666666666, 0.231231, 0.6936929999999999, 0.1091570416994492, 0.1091570416994492, 83.269843363707, 749.4285902733629, 0.0344012356109443, 0.03440123561094431),
    @("MuSCs", "Gnas", "Gcgr", "Resolving-Mac", 3, 1, 360.115397, 1080.346191, 0.3151536087398187, 0.3151536087398188, 1, 0.3333333333333333, 0.1322646666666667, 0.396794, 0.06243808025177026, 0.06243808025177026, 47.63054294573934, 428.674886511654, 0.01967758631413181, 0.01967758631413181),
    @("Resolving-Mac", "Gnas", "Gcgr", "ECs", 3, 1, 76.02905266666666, 228.087158, 0.06653653389046771, 0.06653653389046772, 1, 0.3333333333333333, 0.3424073333333333, 1.027222, 0.1616399685287176, 0.1616399685287176, 26.03290517945289, 234.296146615076, 0.01075496324406515, 0.01075496324406515),
    @("Resolving-Mac", "Gnas", "Gcgr", "FAPs", 3, 1, 76.02905266666666, 228.087158, 0.06653653389046771, 0.06653653389046772, 3, 1, 1.412430333333333, 4.237291, 0.666764909520063, 0.666764909520063, 107.3857402009975, 966.471661808978, 0.0443642259992563, 0.04436422599925631),
    @("Resolving-Mac", "Gnas", "Gcgr", "MuSCs", 3, 1, 76.02905266666666, 228.087158, 0.06653653389046771, 0.06653653389046772, 2, 0.6666666666666666, 0.231231, 0.6936929999999999, 0.1091570416994492, 0.1091570416994492, 17.580273877166, 158.222464894494, 0.007262931204418601, 0.007262931204418603),
    @("Resolving-Mac", "Gnas", "Gcgr", "Resolving-Mac", 3, 1, 76.02905266666666, 228.087158, 0.06653653389046771, 0.06653653389046772, 1, 0.3333333333333333, 0.1322646666666667, 0.396794, 0.06243808025177026, 0.06243808025177026, 10.05595730793911, 90.50361577145199, 0.004154413442727655, 0.004154413442727656),
)

$rowNum = 2
foreach ($row in $data) {
    $colNum = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowNum, $colNum).Value = $val
        $colNum++
    }
    $rowNum++
}